$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O7").Value = 0
$ws.Range("O8").Value = 0
$ws.Range("O11").Value = 0.01579952239990234
$ws.Range("O12").Value = 0
$ws.Range("O14").Value = 0.149911642074585
$ws.Range("O15").Value = 0.001105308532714844
$ws.Range("O17").Value = 0
$ws.Range("O19").Value = 0.005612611770629883
$ws.Range("O20").Value = 0.001584529876708984
$ws.Range("O21").Value = 0
$ws.Range("O25").Value = 0.006570339202880859
$ws.Range("O26").Value = 0.05611276626586914
$ws.Range("O27").Value = 0
$ws.Range("O29").Value = 0
$ws.Range("O30").Value = 0
$ws.Range("O31").Value = 0.004039764404296875
$ws.Range("O33").Value = 0.01628589630126953
$ws.Range("O34").Value = 0.08252167701721191
$ws.Range("O35").Value = 0.03958797454833984
$ws.Range("O36").Value = 0
$ws.Range("O37").Value = 0.3334903717041016
$ws.Range("O38").Value = 0
$ws.Range("O39").Value = 0
$ws.Range("O40").Value = 0
$ws.Range("O41").Value = 0
$ws.Range("O42").Value = 0
$ws.Range("O43").Value = 0
$ws.Range("O44").Value = 0
$ws.Range("O46").Value = 0
$ws.Range("O47").Value = 0.007104396820068359
$ws.Range("O48").Value = 0.001598119735717773
$ws.Range("O49").Value = 0.007414340972900391
$ws.Range("O50").Value = 0
$ws.Range("O52").Value = 0
$ws.Range("O53").Value = 0
$ws.Range("O54").Value = 0
$ws.Range("O55").Value = 0
$ws.Range("O56").Value = 0
$ws.Range("O57").Value = 0
$ws.Range("O59").Value = 0.0101008415222168
$ws.Range("O60").Value = 0
$ws.Range("O61").Value = 0.01000475883483887
$ws.Range("O62").Value = 0
$ws.Range("O63").Value = 0
$ws.Range("O64").Value = 0.01572990417480469
$ws.Range("O65").Value = 0.460355281829834
$ws.Range("O66").Value = 0
$ws.Range("O67").Value = 0
$ws.Range("O68").Value = 0
$ws.Range("O69").Value = 0.01071667671203613
$ws.Range("O71").Value = 1.78434681892395
$ws.Range("O72").Value = 17.63679718971252
$ws.Range("O73").Value = 0.0271143913269043
$ws.Range("O75").Value = 0.001001834869384766
$ws.Range("O77").Value = 0
$ws.Range("O78").Value = 0.01403594017028809
$ws.Range("O79").Value = 0
$ws.Range("O80").Value = 0.001890659332275391
$ws.Range("O81").Value = 0.08059287071228027
$ws.Range("O84").Value = 4.013778209686279
$ws.Range("O85").Value = 0.006505250930786133
$ws.Range("O87").Value = 0.04182934761047363
$ws.Range("O88").Value = 0.04611086845397949
$ws.Range("O89").Value = 0.0009243488311767578
$ws.Range("O90").Value = 0.06563925743103027
$ws.Range("O91").Value = 0.03500270843505859
$ws.Range("O92").Value = 0.02590513229370117
$ws.Range("O93").Value = 0.04880595207214355
$ws.Range("O94").Value = 0.0333554744720459
$ws.Range("O95").Value = 0.0481259822845459
$ws.Range("O96").Value = 0.03818821907043457
$ws.Range("O98").Value = 0
$ws.Range("O99").Value = 0.0270392894744873
$ws.Range("O100").Value = 0
$ws.Range("O101").Value = 0.01790833473205566
$ws.Range("O102").Value = 0.04964923858642578
$ws.Range("O103").Value = 0
$ws.Range("O104").Value = 0.9572534561157227
$ws.Range("O106").Value = 0.02267694473266602
$ws.Range("O108").Value = 0.01648998260498047
$ws.Range("O109").Value = 0.007256031036376953
$ws.Range("O110").Value = 0.01356053352355957
$ws.Range("O111").Value = 0.06697845458984375
$ws.Range("O112").Value = 0.09998917579650879
$ws.Range("O113").Value = 0.09976100921630859
$ws.Range("O114").Value = 0.08468890190124512
$ws.Range("O115").Value = 0
$ws.Range("O116").Value = 0.04265499114990234
$ws.Range("O122").Value = 0.01349592208862305
$ws.Range("O123").Value = 0.00901484489440918
$ws.Range("O129").Value = 0
$ws.Range("O130").Value = 0
$ws.Range("O133").Value = 0
$ws.Range("O136").Value = 0.0009992122650146484
$ws.Range("O138").Value = 0.01255631446838379
$ws.Range("O140").Value = 0.00130915641784668
$ws.Range("O143").Value = 0.004910469055175781
$ws.Range("O144").Value = 0
$ws.Range("O145").Value = 0.001509666442871094
$ws.Range("O146").Value = 0
$ws.Range("O147").Value = 0.001042842864990234
$ws.Range("O148").Value = 0.01930880546569824
$ws.Range("O149").Value = 0
$ws.Range("O150").Value = 0.2493259906768799
$ws.Range("O151").Value = 0.02045822143554688
$ws.Range("O152").Value = 0.009504556655883789
$ws.Range("O162").Value = 0.004445791244506836
$ws.Range("O163").Value = 0.0009999275207519531
$ws.Range("O165").Value = 0.003006935119628906
$ws.Range("O168").Value = 0.008615970611572266
$ws.Range("O170").Value = 0.01304388046264648
$ws.Range("O171").Value = 0.01666164398193359
$ws.Range("O174").Value = 0.01603484153747559
$ws.Range("O176").Value = 0
